$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.953.99'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.09'
$ws.Range("E3").Value = '  +1.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.71'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4575'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3805'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07744'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9765'
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.38'
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.942.20'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.701'
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.960'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06990'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009479'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.69'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.960.01'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.336'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.148.00'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.057'
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.98'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.05'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.604'
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.55'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.840'
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09300'
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8644'
$ws.Range("E32").Value = '  +0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.104'
$ws.Range("E33").Value = '  +0.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.240'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.010'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05685'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.004'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.059'
$ws.Range("E40").Value = '  +11.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.457'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5495'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.328'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002779'
$ws.Range("E45").Value = '  +13.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.173'
$ws.Range("E46").Value = '  +4.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5160'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06936'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.18'
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.59'
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.764'
$ws.Range("E51").Value = '  -0.52%  '
